# Update "paises.xlsx" (Pais sheet) with refreshed COVID-19 statistics and
# re-ordering of a couple of countries whose case counts changed rank order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp on row 1 -----------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 19 de Agosto de 2020 a las 23:06"

# --- Row 4: Estados Unidos --------------------------------------------
$ws.Cells.Item(4, 2).Value = 5687167
$ws.Cells.Item(4, 3).Value = 31193
$ws.Cells.Item(4, 4).Value = 3040396
$ws.Cells.Item(4, 5).Value = 2470844
$ws.Cells.Item(4, 7).Value = 853
$ws.Cells.Item(4, 8).Value = 175927

# --- Row 8: Sudafrica ---------------------------------------------------
$ws.Cells.Item(8, 2).Value = 596060
$ws.Cells.Item(8, 3).Value = 3916
$ws.Cells.Item(8, 4).Value = 491441
$ws.Cells.Item(8, 5).Value = 92196
$ws.Cells.Item(8, 7).Value = 159
$ws.Cells.Item(8, 8).Value = 12423

# --- Row 22: Alemania -----------------------------------------------
$ws.Cells.Item(22, 2).Value = 229688
$ws.Cells.Item(22, 3).Value = 1583
$ws.Cells.Item(22, 5).Value = 16474

# --- Row 32: Israel -------------------------------------------------
$ws.Cells.Item(32, 2).Value = 97969
$ws.Cells.Item(32, 3).Value = 1560
$ws.Cells.Item(32, 4).Value = 73092
$ws.Cells.Item(32, 5).Value = 24096

# --- Row 53: Barein ---------------------------------------------------
$ws.Cells.Item(53, 5).Value = 3482
$ws.Cells.Item(53, 7).Value = 3
$ws.Cells.Item(53, 8).Value = 178

# --- Row 76: Estado de Palestina --------------------------------------
$ws.Cells.Item(76, 5).Value = 7175
$ws.Cells.Item(76, 7).Value = 6
$ws.Cells.Item(76, 8).Value = 119

# --- Row 77: Costa de Marfil ------------------------------------------
$ws.Cells.Item(77, 2).Value = 17232
$ws.Cells.Item(77, 3).Value = 82
$ws.Cells.Item(77, 4).Value = 14422
$ws.Cells.Item(77, 5).Value = 2699
$ws.Cells.Item(77, 7).Value = 1
$ws.Cells.Item(77, 8).Value = 111

# --- Row 93: Guinea -----------------------------------------------------
$ws.Cells.Item(93, 2).Value = 8792
$ws.Cells.Item(93, 3).Value = 77
$ws.Cells.Item(93, 4).Value = 7574
$ws.Cells.Item(93, 5).Value = 1165
$ws.Cells.Item(93, 7).Value = 1
$ws.Cells.Item(93, 8).Value = 53

# --- Row 95: Gabon --------------------------------------------------
$ws.Cells.Item(95, 2).Value = 8319
$ws.Cells.Item(95, 3).Value = 49
$ws.Cells.Item(95, 4).Value = 6614
$ws.Cells.Item(95, 5).Value = 1652

# --- Row 109: Guinea Ecuatorial ---------------------------------------
$ws.Cells.Item(109, 2).Value = 4892
$ws.Cells.Item(109, 3).Value = 71
$ws.Cells.Item(109, 4).Value = 2713
$ws.Cells.Item(109, 5).Value = 2096

# --- Rows 119/120: Somalia & Cabo Verde swapped rank order --------------
# Cabo Verde's case count overtook Somalia's, so the two rows trade places
# (country names swap while keeping the sheet sorted by "Casos totales").
$ws.Cells.Item(119, 1).Value = "Cabo Verde"
$ws.Cells.Item(119, 2).Value = 3321
$ws.Cells.Item(119, 3).Value = 68
$ws.Cells.Item(119, 4).Value = 2442
$ws.Cells.Item(119, 5).Value = 843
$ws.Cells.Item(119, 8).Value = 36

$ws.Cells.Item(120, 1).Value = "Somalia"
$ws.Cells.Item(120, 2).Value = 3265
$ws.Cells.Item(120, 3).Value = 8
$ws.Cells.Item(120, 4).Value = 2396
$ws.Cells.Item(120, 5).Value = 776
$ws.Cells.Item(120, 8).Value = 93

# --- Row 137: Angola --------------------------------------------------
$ws.Cells.Item(137, 2).Value = 2015
$ws.Cells.Item(137, 3).Value = 49
$ws.Cells.Item(137, 4).Value = 698
$ws.Cells.Item(137, 5).Value = 1225
$ws.Cells.Item(137, 7).Value = 2
$ws.Cells.Item(137, 8).Value = 92

# --- Row 138: Sierra Leona --------------------------------------------
$ws.Cells.Item(138, 2).Value = 1961
$ws.Cells.Item(138, 3).Value = 2
$ws.Cells.Item(138, 4).Value = 1531

# --- Row 154: Togo ------------------------------------------------------
$ws.Cells.Item(154, 2).Value = 1190
$ws.Cells.Item(154, 3).Value = 17
$ws.Cells.Item(154, 4).Value = 875
$ws.Cells.Item(154, 5).Value = 288

# --- Row 188: Barbados --------------------------------------------------
$ws.Cells.Item(188, 2).Value = 155
$ws.Cells.Item(188, 3).Value = 2
$ws.Cells.Item(188, 4).Value = 123
$ws.Cells.Item(188, 5).Value = 25

# --- Rows 213/214: Montserrat & Islas Malvinas swapped rank order -------
$ws.Cells.Item(213, 1).Value = "Islas Malvinas"
$ws.Cells.Item(213, 4).Value = 13
$ws.Cells.Item(213, 8).Value = 0

$ws.Cells.Item(214, 1).Value = "Montserrat"
$ws.Cells.Item(214, 4).Value = 12
$ws.Cells.Item(214, 8).Value = 1
